$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four obsolete lookup columns (R:U -> "Dam", "Borne",
# "Quelques pierres", "Quarry"). Deleting the whole columns shifts
# everything to their right (V:AL) left by four, which also updates the
# used-range dimension, the frozen-pane/selection anchors are then fixed
# up explicitly below.
$ws.Range("R:U").Delete() | Out-Null

# Restore the view state captured in the new file.
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("C2").Select() | Out-Null
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("H9").Select() | Out-Null
